$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Mein Verhaltensziel: Adler"
$ws.Range("E6").Value = "https://www.marco-lehmann.de/wp-content/uploads/adler-huhn-1200x675.jpg"
$ws.Range("D6").Value = "WOW"
$ws.Range("C6").Value = "Gestern spürte ich, dass es meinem Körper nicht so gut geht und bin deshalb kurz hingelegen. Und siehe da: Ich schloss die Augen und sah zum ersten Mal so richtig mein Verhaltensziel: der Adler fliegend über mir."

$ws.Range("E11").Select()
